{"js": "// Office.js (Word JavaScript API) script.\n// Rewrites the review's paragraphs to the new review text (different\n// paper, different summary) per the target diff, keeping each\n// paragraph's style (\"Normal\") and position intact.\nconst newTexts = [\n  \"\u26a1\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 08.06.24:\u26a1\ud83d\ude80\",\n  \"Transformers are SSMs: Generalized Models and Efficient Algorithms Through Structured State Space Duality\",\n  \"\u05dc\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05d9\u05e9 \u05e2\u05d5\u05d3 \u05e9\u05dd \u05d5\u05d4\u05d5\u05d0 \ud83e\udd14mamba-2. \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05de\u05ea\u05de\u05e7\u05d3 \u05d1\u05e9\u05db\u05dc\u05d5\u05dc \u05d4\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05e9\u05dc \u05de\u05de\u05d1\u05d4 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9\u05ea \u05e9\u05e2\u05e9\u05ea\u05d4 \u05d4\u05e8\u05d1\u05d4 \u05db\u05d5\u05ea\u05e8\u05d5\u05ea \u05d1\u05d7\u05e6\u05d9 \u05d4\u05e9\u05e0\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05d5\u05d0\u05e0\u05d9 \u05d4\u05e6\u05d8\u05e8\u05e4\u05ea\u05d9 \u05dc\u05d7\u05d2\u05d9\u05d2\u05d4 \u05d5\u05e1\u05e7\u05e8\u05ea\u05d9 \u05d1\u05e2\u05e8\u05da 20 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05d1\u05e0\u05d5\u05e9\u05d0 \u05d4\u05de\u05e8\u05ea\u05e7 \u05d4\u05d6\u05d4.\",\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05e9\u05dc Albert Gu \u05d4\u05ea\u05d5\u05ea\u05d7 \u05de\u05de\u05e9\u05d9\u05da \u05dc\u05d4\u05e2\u05e9\u05d9\u05e8 \u05d0\u05ea \u05e2\u05d5\u05dc\u05dd \u05d4\u05de\u05de\u05d1\u05d4 \u05d5\u05d4\u05e4\u05e2\u05dd \u05d4\u05d5\u05d0 \u05d4\u05d2\u05d9\u05e2 \u05dc\u05db\u05de\u05d4 \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05d5\u05ea. \u05d4\u05d5\u05d0 \u05dc\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 \u05de\u05d2\u05d3\u05d9\u05e8 SSM \u05d1\u05e2\u05dc \u05ea\u05db\u05d5\u05e0\u05d4 N-semi-separable \u05e9\u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05d2\u05d3\u05d9\u05e8 \u05d0\u05ea \u05e6\u05d5\u05e8\u05ea\u05d5 \u05e9\u05dc \u05e7\u05e8\u05e0\u05dc \u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d4 \u05d4\u05de\u05d5\u05e4\u05e2\u05dc \u05e2\u05dc \u05e1\u05d3\u05e8\u05ea \u05d4\u05e7\u05dc\u05d8 \u05d1\u05de\u05d5\u05d3 \u05d4\u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d5\u05e0\u05d9 \u05e9\u05dc SSM (\u05db\u05d0\u05e9\u05e8 \u05de\u05e9\u05ea\u05de\u05e9\u05d9\u05dd \u05d1-SSM \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05de\u05d5\u05e7\u05d1\u05dc). \u05d0\u05dc\u05d7\u05e9 \u05dc\u05db\u05dd \u05d1\u05e1\u05d5\u05d3 \u05e9\u05d1\u05e1\u05d5\u05e4\u05d5 \u05e9\u05dc \u05d3\u05d1\u05e8 \u05d6\u05d4 \u05de\u05ea\u05e0\u05e7\u05d6 \u05dc\u05e6\u05d5\u05e8\u05ea\u05d5 \u05e9\u05dc \u05de\u05d8\u05e8\u05d9\u05e6\u05d4 A. \",\n  \"\u05e9\u05e0\u05d9\u05ea \u05de\u05d0\u05de\u05e8 \u05d7\u05d5\u05e7\u05e8 \u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9 \u05d4-attention \u05d1\u05e4\u05e8\u05d8 \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd \u05dc\u05de\u05e9\u05dc \u05d4\u05e7\u05dc\u05d0\u05e1\u05d9  \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9 ,\u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05dc\u05d0 \u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1, \u05d5\u05e2\u05dd \u05e1\u05d3\u05e8 \u05e9\u05d5\u05e0\u05d4 \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2  \u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d9\u05df \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea Q, K, \u05d5-V. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e4\u05e8\u05e7 \u05d0\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05dc- 3 \u05e9\u05dc\u05d1\u05d9\u05dd \\\"\u05d0\u05d8\u05d5\u05de\u05d9\u05d9\u05dd\\\" (\u05e9\u05db\u05dc \u05d0\u05d7\u05ea \u05de\u05d4\u05dd \u05d4\u05d5\u05d0 \u05de\u05db\u05e4\u05dc\u05d5\u05ea \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea, \u05d0\u05da \u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05de\u05d0\u05d5\u05d3 \u05d2\u05d3\u05d5\u05dc\u05d5\u05ea) \u05d4\u05e9\u05dc\u05d1 \u05d4\u05e9\u05e0\u05d9 \u05d5\u05d4\u05d7\u05e9\u05d5\u05d1 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d4\u05d5\u05d0 \u05de\u05d9\u05e1\u05d5\u05da (masking) \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05ea\u05d0\u05e8 \u05d0\u05d5\u05ea\u05d5 \u05d2\u05dd \u05e2\u05dc\u05d9 \u05d9\u05d3\u05d9 \u05de\u05db\u05e4\u05dc\u05d5\u05ea \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea (Kernel trick). \u05d4\u05de\u05d9\u05e1\u05d5\u05da \u05d4\u05e7\u05d5\u05d6\u05dc\u05d9 (causal) \u05d4\u05d5\u05d0 \u05d7\u05dc\u05e7 \u05de\u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05d4 masking. \u05d4\u05d1\u05d7\u05e0\u05d4 \u05d6\u05d5 \u05d0\u05e4\u05e9\u05e8\u05d4 \u05dc\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05dc\u05d4\u05d5\u05db\u05d9\u05d7 \u05e1\u05d5\u05d2 \u05e9\u05dc \u05e9\u05e7\u05d9\u05dc\u05d5\u05ea \u05d1\u05d9\u05df \u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9 attention \u05de\u05e1\u05d5\u05d9\u05de\u05d9\u05dd \u05dc-SSMs .\",\n  \"\u05d1\u05e0\u05d5\u05e1\u05e3 \u05d4\u05dd \u05de\u05e4\u05ea\u05d7 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d9\u05e2\u05d9\u05dc \u05e9\u05dc \u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d4 \u05d0\u05e8\u05d5\u05db\u05d4 (\u05e9\u05d6\u05d4 \u05d4\u05dc\u05d1 \u05e9\u05dc SSM) \u05d1\u05d7\u05d5\u05de\u05e8\u05d4 \u05e2\u05dd \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d4\u05e0\u05e7\u05e8\u05d0\u05ea 1-semi-separable (\u05e2\u05d1\u05d5\u05e8 \u05de\u05d8\u05e8\u05d9\u05e6\u05d4 A \u05de\u05e6\u05d5\u05e8\u05d4 \u05de\u05e1\u05d5\u05d9\u05de\u05ea).\",\n  \"\u05de\u05d4 \u05d9\u05d5\u05e6\u05d0 \u05dc\u05e0\u05d5 \u05de\u05db\u05dc \u05d4\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4? \u05d4\u05d0\u05e6\u05ea \u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc \u05de\u05de\u05d1\u05d4 (\u05e9\u05d6\u05d4 \u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05de\u05d1\u05d4 2) \u05d5\u05d2\u05dd \u05e4\u05e8\u05d9\u05d9\u05de\u05d5\u05d5\u05e8\u05e7 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9 \u05dc\u05de\u05d9\u05d3\u05d5\u05dc \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05d4\u05e2\u05d5\u05e6\u05de\u05ea\u05d9\u05ea \u05d4\u05d6\u05d5 \u05de\u05e9\u05d5\u05ea\u05e4\u05ea \u05d2\u05dd \u05dc\u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9 \u05d4-attention \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd.\",\n  \"\u05e7\u05e8\u05d9\u05d0\u05d4 \u05de\u05d4\u05e0\u05d4!\",\n  \"https://arxiv.org/abs/2405.21060\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newTexts.length) {\n  throw new Error(\n    \"Expected \" + newTexts.length + \" paragraphs, found \" + paragraphs.items.length\n  );\n}\n\nfor (let i = 0; i < newTexts.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (i === 6) {\n    // This paragraph's run originally carried xml:space=\"preserve\"\n    // (trailing spaces in the old text). Re-using the run via\n    // insertText keeps that flag stuck on the run even though the new\n    // text has no leading/trailing whitespace, so rebuild the\n    // paragraph instead: insert a fresh one with the new text right\n    // after it, then drop the old one.\n    paragraph.insertParagraph(newTexts[i], Word.InsertLocation.after);\n    paragraph.delete();\n  } else {\n    paragraph.insertText(newTexts[i], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# PowerShell (Word COM interop) script.\n# Rewrites the review's paragraphs to the new review text per the\n# target diff, keeping each paragraph (style, order) intact.\n$d = $word.ActiveDocument\n\n$newTexts = @(\n  \"\u26a1\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 08.06.24:\u26a1\ud83d\ude80\",\n  \"Transformers are SSMs: Generalized Models and Efficient Algorithms Through Structured State Space Duality\",\n  \"\u05dc\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05d9\u05e9 \u05e2\u05d5\u05d3 \u05e9\u05dd \u05d5\u05d4\u05d5\u05d0 \ud83e\udd14mamba-2. \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05de\u05ea\u05de\u05e7\u05d3 \u05d1\u05e9\u05db\u05dc\u05d5\u05dc \u05d4\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05e9\u05dc \u05de\u05de\u05d1\u05d4 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9\u05ea \u05e9\u05e2\u05e9\u05ea\u05d4 \u05d4\u05e8\u05d1\u05d4 \u05db\u05d5\u05ea\u05e8\u05d5\u05ea \u05d1\u05d7\u05e6\u05d9 \u05d4\u05e9\u05e0\u05d4 \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05d5\u05d0\u05e0\u05d9 \u05d4\u05e6\u05d8\u05e8\u05e4\u05ea\u05d9 \u05dc\u05d7\u05d2\u05d9\u05d2\u05d4 \u05d5\u05e1\u05e7\u05e8\u05ea\u05d9 \u05d1\u05e2\u05e8\u05da 20 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05d1\u05e0\u05d5\u05e9\u05d0 \u05d4\u05de\u05e8\u05ea\u05e7 \u05d4\u05d6\u05d4.\",\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05e9\u05dc Albert Gu \u05d4\u05ea\u05d5\u05ea\u05d7 \u05de\u05de\u05e9\u05d9\u05da \u05dc\u05d4\u05e2\u05e9\u05d9\u05e8 \u05d0\u05ea \u05e2\u05d5\u05dc\u05dd \u05d4\u05de\u05de\u05d1\u05d4 \u05d5\u05d4\u05e4\u05e2\u05dd \u05d4\u05d5\u05d0 \u05d4\u05d2\u05d9\u05e2 \u05dc\u05db\u05de\u05d4 \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05e0\u05d5\u05ea. \u05d4\u05d5\u05d0 \u05dc\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 \u05de\u05d2\u05d3\u05d9\u05e8 SSM \u05d1\u05e2\u05dc \u05ea\u05db\u05d5\u05e0\u05d4 N-semi-separable \u05e9\u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05d2\u05d3\u05d9\u05e8 \u05d0\u05ea \u05e6\u05d5\u05e8\u05ea\u05d5 \u05e9\u05dc \u05e7\u05e8\u05e0\u05dc \u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d4 \u05d4\u05de\u05d5\u05e4\u05e2\u05dc \u05e2\u05dc \u05e1\u05d3\u05e8\u05ea \u05d4\u05e7\u05dc\u05d8 \u05d1\u05de\u05d5\u05d3 \u05d4\u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d5\u05e0\u05d9 \u05e9\u05dc SSM (\u05db\u05d0\u05e9\u05e8 \u05de\u05e9\u05ea\u05de\u05e9\u05d9\u05dd \u05d1-SSM \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05de\u05d5\u05e7\u05d1\u05dc). \u05d0\u05dc\u05d7\u05e9 \u05dc\u05db\u05dd \u05d1\u05e1\u05d5\u05d3 \u05e9\u05d1\u05e1\u05d5\u05e4\u05d5 \u05e9\u05dc \u05d3\u05d1\u05e8 \u05d6\u05d4 \u05de\u05ea\u05e0\u05e7\u05d6 \u05dc\u05e6\u05d5\u05e8\u05ea\u05d5 \u05e9\u05dc \u05de\u05d8\u05e8\u05d9\u05e6\u05d4 A. \",\n  \"\u05e9\u05e0\u05d9\u05ea \u05de\u05d0\u05de\u05e8 \u05d7\u05d5\u05e7\u05e8 \u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9 \u05d4-attention \u05d1\u05e4\u05e8\u05d8 \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd \u05dc\u05de\u05e9\u05dc \u05d4\u05e7\u05dc\u05d0\u05e1\u05d9  \u05d4\u05dc\u05d9\u05e0\u05d0\u05e8\u05d9 ,\u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05dc\u05d0 \u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1, \u05d5\u05e2\u05dd \u05e1\u05d3\u05e8 \u05e9\u05d5\u05e0\u05d4 \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2  \u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d9\u05df \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea Q, K, \u05d5-V. \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e4\u05e8\u05e7 \u05d0\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1 \u05dc- 3 \u05e9\u05dc\u05d1\u05d9\u05dd \"\"\u05d0\u05d8\u05d5\u05de\u05d9\u05d9\u05dd\"\" (\u05e9\u05db\u05dc \u05d0\u05d7\u05ea \u05de\u05d4\u05dd \u05d4\u05d5\u05d0 \u05de\u05db\u05e4\u05dc\u05d5\u05ea \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea, \u05d0\u05da \u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05de\u05d0\u05d5\u05d3 \u05d2\u05d3\u05d5\u05dc\u05d5\u05ea) \u05d4\u05e9\u05dc\u05d1 \u05d4\u05e9\u05e0\u05d9 \u05d5\u05d4\u05d7\u05e9\u05d5\u05d1 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05d4\u05d5\u05d0 \u05de\u05d9\u05e1\u05d5\u05da (masking) \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05ea\u05d0\u05e8 \u05d0\u05d5\u05ea\u05d5 \u05d2\u05dd \u05e2\u05dc\u05d9 \u05d9\u05d3\u05d9 \u05de\u05db\u05e4\u05dc\u05d5\u05ea \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea (Kernel trick). \u05d4\u05de\u05d9\u05e1\u05d5\u05da \u05d4\u05e7\u05d5\u05d6\u05dc\u05d9 (causal) \u05d4\u05d5\u05d0 \u05d7\u05dc\u05e7 \u05de\u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05d4 masking. \u05d4\u05d1\u05d7\u05e0\u05d4 \u05d6\u05d5 \u05d0\u05e4\u05e9\u05e8\u05d4 \u05dc\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05dc\u05d4\u05d5\u05db\u05d9\u05d7 \u05e1\u05d5\u05d2 \u05e9\u05dc \u05e9\u05e7\u05d9\u05dc\u05d5\u05ea \u05d1\u05d9\u05df \u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9 attention \u05de\u05e1\u05d5\u05d9\u05de\u05d9\u05dd \u05dc-SSMs .\",\n  \"\u05d1\u05e0\u05d5\u05e1\u05e3 \u05d4\u05dd \u05de\u05e4\u05ea\u05d7 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d9\u05e2\u05d9\u05dc \u05e9\u05dc \u05e7\u05d5\u05e0\u05d1\u05d5\u05dc\u05d5\u05e6\u05d9\u05d4 \u05d0\u05e8\u05d5\u05db\u05d4 (\u05e9\u05d6\u05d4 \u05d4\u05dc\u05d1 \u05e9\u05dc SSM) \u05d1\u05d7\u05d5\u05de\u05e8\u05d4 \u05e2\u05dd \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d4\u05e0\u05e7\u05e8\u05d0\u05ea 1-semi-separable (\u05e2\u05d1\u05d5\u05e8 \u05de\u05d8\u05e8\u05d9\u05e6\u05d4 A \u05de\u05e6\u05d5\u05e8\u05d4 \u05de\u05e1\u05d5\u05d9\u05de\u05ea).\",\n  \"\u05de\u05d4 \u05d9\u05d5\u05e6\u05d0 \u05dc\u05e0\u05d5 \u05de\u05db\u05dc \u05d4\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4? \u05d4\u05d0\u05e6\u05ea \u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc \u05de\u05de\u05d1\u05d4 (\u05e9\u05d6\u05d4 \u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05de\u05d1\u05d4 2) \u05d5\u05d2\u05dd \u05e4\u05e8\u05d9\u05d9\u05de\u05d5\u05d5\u05e8\u05e7 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9 \u05dc\u05de\u05d9\u05d3\u05d5\u05dc \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4 \u05d4\u05e2\u05d5\u05e6\u05de\u05ea\u05d9\u05ea \u05d4\u05d6\u05d5 \u05de\u05e9\u05d5\u05ea\u05e4\u05ea \u05d2\u05dd \u05dc\u05de\u05e0\u05d2\u05e0\u05d5\u05e0\u05d9 \u05d4-attention \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd.\",\n  \"\u05e7\u05e8\u05d9\u05d0\u05d4 \u05de\u05d4\u05e0\u05d4!\",\n  \"https://arxiv.org/abs/2405.21060\"\n)\n\nfor ($i = 0; $i -lt $newTexts.Count; $i++) {\n  $p = $d.Paragraphs.Item($i + 1)\n  $rng = $p.Range\n  # Exclude the trailing paragraph mark so only the visible text is\n  # replaced; this also avoids carrying over a stale\n  # xml:space=\"preserve\" flag from the old run onto the new text.\n  $rng.End = $rng.End - 1\n  $rng.Text = $newTexts[$i]\n}\n"}
